# Atualizei dados da bibi - faturamento anual 2025 (linha 9)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3091295.22
$ws.Range("C9").Value = 485916.99
$ws.Range("D9").Value = 3577212.21
$ws.Range("E9").Value = 13.58367805638235
$ws.Range("F9").Value = 86.41632194361766
$ws.Range("G9").Value = -53.03854581031157
$ws.Range("H9").Value = -44.1754812375133
$ws.Range("I9").Value = 30949
$ws.Range("J9").Value = 1312
$ws.Range("K9").Value = 32261
$ws.Range("L9").Value = 22260
$ws.Range("M9").Value = 160.7013571428571
$ws.Range("N9").Value = 9.714089643186652
